# Fix bug: biyi kouyi use the same pic.
#
# The whole "process" diagram on slide 2 (all folded-corner boxes and
# their connector arrows) was nudged by a constant offset
# (+100667 EMU horizontally, +109057 EMU vertically) so it lines up
# with the picture behind it. Every shape keeps its original size -
# only its position (Left/Top) changes.
#
# Target absolute positions below are the exact EMU values from the
# target OOXML, expressed in points (EMU / 12700). The COM layer here
# stores Shape.Left/Top with single (32-bit) float precision and
# truncates when converting back to EMU, so the literals are tuned
# (nearest float32 that truncates to the right EMU) rather than the
# plain EMU/12700 ratio, to land exactly on the target integers.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$targets = @{
    5  = @(24.2850399017334,   71.3125991821289)
    6  = @(275.13726806640625, 64.60952758789062)
    7  = @(523.3289184570312,  71.25543975830078)
    9  = @(24.526063919067383, 175.07362365722656)
    10 = @(207.0072479248047,  175.07362365722656)
    11 = @(340.8478088378906,  175.07362365722656)
    12 = @(523.3289794921875,  175.07362365722656)
    13 = @(24.2850399017334,   272.2922058105469)
    19 = @(196.04969787597656, 272.2922058105469)
    20 = @(399.7049865722656,  272.4530029296875)
    21 = @(603.3601684570312,  272.2922058105469)
    25 = @(135.76315307617188, 87.21968841552734)
    28 = @(406.55828857421875, 87.16252136230469)
    37 = @(286.1728515625,     190.980712890625)
    39 = @(468.4822082519531,  190.980712890625)
    40 = @(71.13819122314453,  288.1993103027344)
    47 = @(292.8049011230469,  -38.20543670654297)
    50 = @(275.2153015136719,  288.1993103027344)
    52 = @(478.8705749511719,  288.1993103027344)
    56 = @(309.82080078125,    -118.407958984375)
    75 = @(152.1604766845703,  190.980712890625)
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    $key = $sh.Id
    if ($targets.ContainsKey($key)) {
        $pos = $targets[$key]
        $sh.Left = $pos[0]
        $sh.Top = $pos[1]
    }
}
